$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell -> new value
$updates = @{
    "segundas" = @{
        "B2" = 0.02179753267665111
        "D2" = 0.01879642144659549
        "F2" = 0.08318987380918343
        "B3" = 0.01633651666480997
        "D3" = 0.01417228597347231
        "F3" = 0.07515066141732571
    }
    "tercas" = @{
        "B2" = 0.03031291637210963
        "D2" = 0.04799263465484215
        "B3" = 0.02584053947893853
        "D3" = 0.04590921430573769
    }
    "quartas" = @{
        "B2" = 0.0291473263824146
        "B3" = 0.02577321002307439
    }
    "quintas" = @{
        "B2" = 0.01903664042968252
        "D2" = 0.02641991141741556
        "B3" = 0.01858664733778813
        "D3" = 0.02629955037142393
    }
    "sextas" = @{
        "B2" = 0.02601523430715038
        "D2" = 0.02017375267478982
        "F2" = 0.03241334987898496
        "B3" = 0.02594801522922061
        "D3" = 0.01947339507150611
        "F3" = 0.02393297254669709
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
